$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.942.71"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.833.71"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'244.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'0.6940"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'0.9985"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.07688"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "'0.3050"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'23.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "'0.07811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'92.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "1.832.71"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").Value = "'5.093"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "'0.6862"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'6.495"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "'0.000008254"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "28.940.58"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "'242.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").Value = "2.073.46"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "'12.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'7.481"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "'0.9989"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.1493"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Value = "'158.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "'8.758"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").Value = "'18.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "'4.226"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'4.153"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'0.05121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "'0.7734"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "'1.854"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'1.142"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").Value = "'2.692"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.275.35"
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01864"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'2.708"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").Value = "'0.9514"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.54%  "
$ws.Range("D42").Value = "'6.148"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("D43").Value = "'106.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "'0.9984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'9.668"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000123"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5170"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.973.83"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'63.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.00%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.751"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'6.972"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
